# Update crypto price/volume data per the Thu Jul 11 23:49:41 UTC 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell holding the default (unmodified) style used to restore styling
# after temporarily forcing text format on numeric-looking price strings.
$defaultStyleCell = $ws.Range("D4")

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "57.324.44"
$c.Style = $defaultStyleCell.Style
$ws.Range("E2").Value = "  -0.84%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.097.47"
$c.Style = $defaultStyleCell.Style
$ws.Range("E3").Value = "  -0.22%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "524.51"
$c.Style = $defaultStyleCell.Style
$ws.Range("E5").Value = "  -0.07%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "135.90"
$c.Style = $defaultStyleCell.Style
$ws.Range("E6").Value = "  -4.38%  "

$ws.Range("E7").Value = "  -0.01%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "3.095.55"
$c.Style = $defaultStyleCell.Style
$ws.Range("E8").Value = "  -0.32%  "

$ws.Range("E9").Value = "  +2.12%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "7.31"
$c.Style = $defaultStyleCell.Style
$ws.Range("E10").Value = "  +1.14%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.107"
$c.Style = $defaultStyleCell.Style
$ws.Range("E11").Value = "  -1.37%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.395"
$c.Style = $defaultStyleCell.Style
$ws.Range("E12").Value = "  +1.47%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "3.631.54"
$c.Style = $defaultStyleCell.Style
$ws.Range("E13").Value = "  -0.23%  "

$ws.Range("E14").Value = "  +2.23%  "

$ws.Range("E15").Value = "  -2.21%  "

$ws.Range("E16").Value = "  -1.30%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "57.432.05"
$c.Style = $defaultStyleCell.Style
$ws.Range("E17").Value = "  -0.80%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.081.11"
$c.Style = $defaultStyleCell.Style
$ws.Range("E18").Value = "  -0.77%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "5.90"
$c.Style = $defaultStyleCell.Style
$ws.Range("E19").Value = "  -3.36%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "12.36"
$c.Style = $defaultStyleCell.Style
$ws.Range("E20").Value = "  -3.51%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "7.83"
$c.Style = $defaultStyleCell.Style
$ws.Range("E21").Value = "  -2.70%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "346.89"
$c.Style = $defaultStyleCell.Style
$ws.Range("E22").Value = "  +1.63%  "

$ws.Range("E23").Value = "  +0.07%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "67.51"
$c.Style = $defaultStyleCell.Style
$ws.Range("E24").Value = "  +0.81%  "

$ws.Range("E25").Value = "  -2.94%  "

$ws.Range("E26").Value = "  -2.42%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.01"
$c.Style = $defaultStyleCell.Style
$ws.Range("E27").Value = "  +0.91%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.0₃0887"
$c.Style = $defaultStyleCell.Style
$ws.Range("E28").Value = "  -3.78%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = $defaultStyleCell.Style
$ws.Range("E29").Value = "  -0.10%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "7.35"
$c.Style = $defaultStyleCell.Style
$ws.Range("E30").Value = "  +2.00%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.87"
$c.Style = $defaultStyleCell.Style
$ws.Range("E31").Value = "  -0.07%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "6.01"
$c.Style = $defaultStyleCell.Style
$ws.Range("E32").Value = "  -7.62%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "20.63"
$c.Style = $defaultStyleCell.Style
$ws.Range("E33").Value = "  -1.93%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.93"
$c.Style = $defaultStyleCell.Style
$ws.Range("E34").Value = "  +6.38%  "

$ws.Range("E35").Value = "  -4.24%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "158.55"
$c.Style = $defaultStyleCell.Style
$ws.Range("E36").Value = "  +1.58%  "

$ws.Range("E37").Value = "  -1.95%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "25.78"
$c.Style = $defaultStyleCell.Style
$ws.Range("E38").Value = "  -5.28%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.23"
$c.Style = $defaultStyleCell.Style
$ws.Range("E39").Value = "  -2.21%  "

$ws.Range("E40").Value = "  +5.42%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.0657"
$c.Style = $defaultStyleCell.Style
$ws.Range("E41").Value = "  -1.14%  "

$ws.Range("E42").Value = "  +2.37%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.696"
$c.Style = $defaultStyleCell.Style
$ws.Range("E43").Value = "  +1.71%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.374.03"
$c.Style = $defaultStyleCell.Style
$ws.Range("E44").Value = "  +3.36%  "

$ws.Range("E45").Value = "  -0.92%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = $defaultStyleCell.Style

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0266"
$c.Style = $defaultStyleCell.Style
$ws.Range("E47").Value = "  +1.82%  "

$ws.Range("E48").Value = "  -1.95%  "

$ws.Range("E49").Value = "  -1.62%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "19.63"
$c.Style = $defaultStyleCell.Style
$ws.Range("E50").Value = "  -4.67%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.755"
$c.Style = $defaultStyleCell.Style
$ws.Range("E51").Value = "  +2.36%  "
